$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3256.5557
$ws.Range("I76").Value = 3187.6875
$ws.Range("K76").Value = 3187.6875
$ws.Range("M76").Value = -2872.6875
$ws.Range("H79").Value = 3256.5557
$ws.Range("I79").Value = 3187.6875
$ws.Range("K79").Value = 3187.6875
$ws.Range("M79").Value = -2095.6875
$ws.Range("H112").Value = 6945789.5
$ws.Range("J112").Value = 8065639.5
$ws.Range("L112").Value = 24196918.5
$ws.Range("N112").Value = -24199134.5
$ws.Range("H137").Value = 4692.381
$ws.Range("I137").Value = 6090.5835
$ws.Range("K137").Value = 18271.7505
$ws.Range("M137").Value = -15721.7505
$ws.Range("H141").Value = 539734.75
$ws.Range("I141").Value = 1947.3334
$ws.Range("J141").Value = 2153097
$ws.Range("K141").Value = 5842.0002
$ws.Range("L141").Value = 6459291
$ws.Range("M141").Value = -662.0002000000004
$ws.Range("N141").Value = -6469651

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6371.371
$ws.Range("I32").Value = 4358.5283
$ws.Range("J32").Value = 18224.777
$ws.Range("K32").Value = 4358.5283
$ws.Range("L32").Value = 18224.777
$ws.Range("M32").Value = -4071.5283
$ws.Range("N32").Value = -18798.777
$ws.Range("H63").Value = 2673.5833
$ws.Range("I63").Value = 2367
$ws.Range("J63").Value = 3593.3333
$ws.Range("K63").Value = 2367
$ws.Range("L63").Value = 3593.3333
$ws.Range("M63").Value = -1681
$ws.Range("N63").Value = -4965.3333
$ws.Range("H66").Value = 2673.5833
$ws.Range("I66").Value = 2367
$ws.Range("J66").Value = 3593.3333
$ws.Range("K66").Value = 11835
$ws.Range("L66").Value = 17966.6665
$ws.Range("M66").Value = -8403
$ws.Range("N66").Value = -24830.6665
$ws.Range("H74").Value = 1309.3334
$ws.Range("I74").Value = 1130.6666
$ws.Range("J74").Value = 1666.6666
$ws.Range("K74").Value = 1130.6666
$ws.Range("L74").Value = 1666.6666
$ws.Range("M74").Value = -256.6666
$ws.Range("N74").Value = -3414.6666
$ws.Range("H77").Value = 1309.3334
$ws.Range("I77").Value = 1130.6666
$ws.Range("J77").Value = 1666.6666
$ws.Range("K77").Value = 5653.333000000001
$ws.Range("L77").Value = 8333.333000000001
$ws.Range("M77").Value = -1285.333000000001
$ws.Range("N77").Value = -17069.333
$ws.Range("H94").Value = 26862.154
$ws.Range("J94").Value = 26862.154
$ws.Range("L94").Value = 26862.154
$ws.Range("N94").Value = -28664.154
$ws.Range("H132").Value = 37042388
$ws.Range("I132").Value = 47624364
$ws.Range("J132").Value = 5465.6665
$ws.Range("K132").Value = 142873092
$ws.Range("L132").Value = 16396.9995
$ws.Range("M132").Value = -142870562
$ws.Range("N132").Value = -21456.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3427.7222
$ws.Range("I134").Value = 2666.5557
$ws.Range("K134").Value = 7999.6671
$ws.Range("M134").Value = -5464.6671
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 30854
$ws.Range("J137").Value = 29999
$ws.Range("L137").Value = 29999
$ws.Range("N137").Value = -40199
$ws.Range("H138").Value = 28000
$ws.Range("J138").Value = 28000
$ws.Range("L138").Value = 28000
$ws.Range("N138").Value = -38280
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3339.3333
$ws.Range("I31").Value = 1801.862
$ws.Range("J31").Value = 7798
$ws.Range("K31").Value = 1801.862
$ws.Range("L31").Value = 7798
$ws.Range("M31").Value = -1506.862
$ws.Range("N31").Value = -8388
$ws.Range("H34").Value = 3339.3333
$ws.Range("I34").Value = 1801.862
$ws.Range("J34").Value = 7798
$ws.Range("K34").Value = 1801.862
$ws.Range("L34").Value = 7798
$ws.Range("M34").Value = -1599.862
$ws.Range("N34").Value = -8202
$ws.Range("H134").Value = 2198.1428
$ws.Range("I134").Value = 1394.3103
$ws.Range("J134").Value = 6083.3335
$ws.Range("K134").Value = 4182.9309
$ws.Range("L134").Value = 18250.0005
$ws.Range("M134").Value = -1647.9309
$ws.Range("N134").Value = -23320.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1233.8
$ws.Range("I5").Value = 697.4091
$ws.Range("J5").Value = 2141.5386
$ws.Range("K5").Value = 2092.2273
$ws.Range("L5").Value = 6424.6158
$ws.Range("M5").Value = -1980.2273
$ws.Range("N5").Value = -6648.6158
$ws.Range("H122").Value = 1243.25
$ws.Range("I122").Value = 746.25
$ws.Range("J122").Value = 1408.9166
$ws.Range("K122").Value = 6716.25
$ws.Range("L122").Value = 12680.2494
$ws.Range("M122").Value = -4266.25
$ws.Range("N122").Value = -17580.2494
$ws.Range("H131").Value = 1101.6129
$ws.Range("I131").Value = 883.46155
$ws.Range("J131").Value = 2236
$ws.Range("K131").Value = 2650.38465
$ws.Range("L131").Value = 6708
$ws.Range("M131").Value = 2389.61535
$ws.Range("N131").Value = -16788
$ws.Range("H132").Value = 2321.35
$ws.Range("I132").Value = 1788.4667
$ws.Range("J132").Value = 3920
$ws.Range("K132").Value = 16096.2003
$ws.Range("L132").Value = 35280
$ws.Range("M132").Value = -13566.2003
$ws.Range("N132").Value = -40340
$ws.Range("H135").Value = 1233.8
$ws.Range("I135").Value = 697.4091
$ws.Range("J135").Value = 2141.5386
$ws.Range("K135").Value = 6276.6819
$ws.Range("L135").Value = 19273.8474
$ws.Range("M135").Value = -3741.6819
$ws.Range("N135").Value = -24343.8474
$ws.Range("H140").Value = 13890889
$ws.Range("I140").Value = 27778584
$ws.Range("J140").Value = 3193.6667
$ws.Range("K140").Value = 83335752
$ws.Range("L140").Value = 9581.000100000001
$ws.Range("M140").Value = -83330572
$ws.Range("N140").Value = -19941.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6549.357
$ws.Range("I70").Value = 7208.273
$ws.Range("J70").Value = 4133.3335
$ws.Range("K70").Value = 7208.273
$ws.Range("L70").Value = 4133.3335
$ws.Range("M70").Value = -6938.273
$ws.Range("N70").Value = -4673.3335
$ws.Range("H73").Value = 6549.357
$ws.Range("I73").Value = 7208.273
$ws.Range("J73").Value = 4133.3335
$ws.Range("K73").Value = 7208.273
$ws.Range("L73").Value = 4133.3335
$ws.Range("M73").Value = -6272.273
$ws.Range("N73").Value = -6005.3335
$ws.Range("H132").Value = 3060.3809
$ws.Range("I132").Value = 2304.6667
$ws.Range("J132").Value = 4949.6665
$ws.Range("K132").Value = 6914.000100000001
$ws.Range("L132").Value = 14848.9995
$ws.Range("M132").Value = -4384.000100000001
$ws.Range("N132").Value = -19908.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2820.05
$ws.Range("I132").Value = 2217.1667
$ws.Range("J132").Value = 3078.4285
$ws.Range("K132").Value = 6651.500100000001
$ws.Range("L132").Value = 9235.2855
$ws.Range("M132").Value = -4121.500100000001
$ws.Range("N132").Value = -14295.2855
$ws.Range("H133").Value = 29500
$ws.Range("J133").Value = 29500
$ws.Range("L133").Value = 29500
$ws.Range("N133").Value = -34560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8739
$ws.Range("I132").Value = 2409.4546
$ws.Range("J132").Value = 18685.428
$ws.Range("K132").Value = 7228.3638
$ws.Range("L132").Value = 56056.284
$ws.Range("M132").Value = -4698.3638
$ws.Range("N132").Value = -61116.284
$ws.Range("H138").Value = 29495
$ws.Range("J138").Value = 29495
$ws.Range("L138").Value = 29495
$ws.Range("N138").Value = -39775
